$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a text (string) value even when the
    # content looks like a number (e.g. "596.54"), matching the
    # original inline/shared-string cell type, then restore the
    # default "Normal" style so no stray formatting is introduced.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "71.232.46"
$ws.Range("E2").Value = "  +6.77%  "
$ws.Range("D3").Value = "3.677.92"
$ws.Range("E3").Value = "  +18.78%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue $ws.Range("D5") "596.54"
$ws.Range("E5").Value = "  +3.40%  "
Set-TextValue $ws.Range("D6") "183.41"
$ws.Range("E6").Value = "  +5.98%  "
$ws.Range("D7").Value = "3.677.38"
$ws.Range("E7").Value = "  +18.82%  "
$ws.Range("E8").Value = "  +0.03%  "
Set-TextValue $ws.Range("D9") "0.535"
$ws.Range("E9").Value = "  +4.28%  "
$ws.Range("E10").Value = "  +7.55%  "
Set-TextValue $ws.Range("D11") "6.57"
$ws.Range("E11").Value = "  +3.65%  "
$ws.Range("E12").Value = "  +6.50%  "
Set-TextValue $ws.Range("D13") "39.73"
$ws.Range("E13").Value = "  +10.19%  "
$ws.Range("E14").Value = "  +5.88%  "
$ws.Range("D15").Value = "4.291.48"
$ws.Range("E15").Value = "  +18.86%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "71.252.96"
$ws.Range("E16").Value = "  +6.91%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.676.40"
$ws.Range("E17").Value = "  +18.83%  "
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("E19").Value = "  +7.47%  "
Set-TextValue $ws.Range("D20") "16.91"
$ws.Range("E20").Value = "  -0.02%  "
Set-TextValue $ws.Range("D21") "517.65"
$ws.Range("E21").Value = "  +6.06%  "
Set-TextValue $ws.Range("D22") "9.17"
$ws.Range("E22").Value = "  +17.62%  "
$ws.Range("E23").Value = "  +8.11%  "
Set-TextValue $ws.Range("D24") "87.69"
$ws.Range("E24").Value = "  +5.18%  "
Set-TextValue $ws.Range("D25") "13.52"
$ws.Range("E25").Value = "  +6.58%  "
$ws.Range("E26").Value = "  +7.70%  "
Set-TextValue $ws.Range("D27") "10.82"
$ws.Range("E27").Value = "  +7.63%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +12.24%  "
Set-TextValue $ws.Range("D30") "8.13"
$ws.Range("E30").Value = "  +2.38%  "
Set-TextValue $ws.Range("D31") "31.81"
$ws.Range("E31").Value = "  +13.76%  "
$ws.Range("E32").Value = "  +6.37%  "
$ws.Range("E33").Value = "  +17.23%  "
$ws.Range("E34").Value = "  +3.58%  "
$ws.Range("E35").Value = "  +0.14%  "
Set-TextValue $ws.Range("D36") "6.17"
$ws.Range("E36").Value = "  +10.50%  "
$ws.Range("E37").Value = "  +7.96%  "
$ws.Range("E38").Value = "  +10.89%  "
Set-TextValue $ws.Range("D39") "2.14"
$ws.Range("E39").Value = "  +9.24%  "
Set-TextValue $ws.Range("D40") "50.77"
$ws.Range("E40").Value = "  +3.33%  "
Set-TextValue $ws.Range("D41") "46.05"
$ws.Range("E41").Value = "  -6.07%  "
$ws.Range("E42").Value = "  +4.05%  "
$ws.Range("D43").Value = "3.173.78"
$ws.Range("E43").Value = "  +14.25%  "
Set-TextValue $ws.Range("D44") "8.79"
$ws.Range("E44").Value = "  +6.45%  "
$ws.Range("E45").Value = "  +6.26%  "
Set-TextValue $ws.Range("D46") "400.14"
$ws.Range("E46").Value = "  +8.59%  "
Set-TextValue $ws.Range("D47") "0.0367"
$ws.Range("E47").Value = "  +6.48%  "
Set-TextValue $ws.Range("D48") "28.05"
$ws.Range("E48").Value = "  +14.65%  "
Set-TextValue $ws.Range("D49") "136.06"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("E50").Value = "  -0.01%  "
Set-TextValue $ws.Range("D51") "2.44"
$ws.Range("E51").Value = "  +11.81%  "
